{"js": "// Sequential, in-document-order replacement of the date line and every\n// division-problem cell's text, matching the target edit exactly\n// (several source strings, e.g. \"22\u00f72=\", repeat but map to different\n// targets depending on position, so we walk paragraphs in document\n// order rather than doing a global find/replace).\nconst replacements = [\n  \"2023-10-25 Wednesday\",\n  \"34\u00f74=\",\n  \"68\u00f75=\",\n  \"67\u00f74=\",\n  \"24\u00f76=\",\n  \"21\u00f72=\",\n  \"92\u00f77=\",\n  \"28\u00f79=\",\n  \"54\u00f77=\",\n  \"55\u00f78=\",\n  \"57\u00f79=\",\n  \"39\u00f72=\",\n  \"89\u00f74=\",\n  \"17\u00f78=\",\n  \"29\u00f74=\",\n  \"54\u00f73=\",\n  \"27\u00f78=\",\n  \"48\u00f77=\",\n  \"63\u00f72=\",\n  \"73\u00f77=\",\n  \"28\u00f73=\",\n  \"68\u00f79=\",\n  \"28\u00f77=\",\n  \"85\u00f78=\",\n  \"35\u00f73=\",\n  \"74\u00f79=\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet i = 0;\nfor (const p of paragraphs.items) {\n  if (p.text && p.text.trim().length > 0) {\n    if (i >= replacements.length) {\n      break;\n    }\n    p.insertText(replacements[i], \"Replace\");\n    i++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Sequential, in-document-order replacement of the date line and every\n# division-problem cell's text, matching the target edit exactly.\n# Several source strings (e.g. \"22\u00f72=\") repeat but map to different\n# targets depending on position, so each pair is located with Find\n# starting just after the end of the previous replacement, walking\n# forward through the document instead of doing a single global\n# Find/Replace (which could not disambiguate duplicate source text).\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2023-10-24 Tuesday\", \"2023-10-25 Wednesday\"),\n    @(\"56\u00f74=\", \"34\u00f74=\"),\n    @(\"16\u00f76=\", \"68\u00f75=\"),\n    @(\"14\u00f75=\", \"67\u00f74=\"),\n    @(\"22\u00f72=\", \"24\u00f76=\"),\n    @(\"22\u00f72=\", \"21\u00f72=\"),\n    @(\"98\u00f77=\", \"92\u00f77=\"),\n    @(\"78\u00f75=\", \"28\u00f79=\"),\n    @(\"51\u00f79=\", \"54\u00f77=\"),\n    @(\"50\u00f76=\", \"55\u00f78=\"),\n    @(\"75\u00f74=\", \"57\u00f79=\"),\n    @(\"87\u00f78=\", \"39\u00f72=\"),\n    @(\"23\u00f77=\", \"89\u00f74=\"),\n    @(\"84\u00f73=\", \"17\u00f78=\"),\n    @(\"58\u00f76=\", \"29\u00f74=\"),\n    @(\"47\u00f76=\", \"54\u00f73=\"),\n    @(\"54\u00f79=\", \"27\u00f78=\"),\n    @(\"11\u00f74=\", \"48\u00f77=\"),\n    @(\"42\u00f79=\", \"63\u00f72=\"),\n    @(\"38\u00f77=\", \"73\u00f77=\"),\n    @(\"71\u00f72=\", \"28\u00f73=\"),\n    @(\"66\u00f73=\", \"68\u00f79=\"),\n    @(\"79\u00f75=\", \"28\u00f77=\"),\n    @(\"10\u00f78=\", \"85\u00f78=\"),\n    @(\"48\u00f75=\", \"35\u00f73=\"),\n    @(\"69\u00f79=\", \"74\u00f79=\")\n)\n\n$cursor = 0\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $rng = $d.Range($cursor, $d.Content.End)\n    $rng.Find.ClearFormatting()\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 0\n    $rng.Find.Text = $old\n\n    if ($rng.Find.Execute()) {\n        $rng.Text = $new\n        $cursor = $rng.End\n    }\n}\n"}
